$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The columns "Estacion mas cercana 6" and "Estacion mas cercana 7" (currently
# columns Q:R) need to move so they sit right after "Estacion mas cercana 5"
# (column K), i.e. become the new columns L:M. The "Inicio estacion mas
# cercana 1..5" columns (currently L:P) then shift right to become N:R.
#
# Insert two blank columns before L, which pushes the existing L:T block
# (Inicio 1-5, Estacion 6-7, Inicio 6-7) to N:V. The data we want at L:M
# (old Estacion 6/7) is now sitting at S:T, so copy it into place and then
# remove the now-duplicated S:T columns.

$ws.Columns("L:M").Insert(-4161)
$ws.Range("S1:T59").Copy()
$ws.Range("L1").PasteSpecial(-4142)
$ws.Columns("S:T").Delete(-4161)
